$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.399.35"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.849.19"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.29"
$ws.Range("E5").Value = "  +1.82%  "

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").Value = "  +2.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2747"
$ws.Range("E8").Value = "  +2.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06312"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.56"
$ws.Range("E10").Value = "  +10.31%  "

$ws.Range("D11").Value = "1.827.63"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07465"
$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.942"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.35"
$ws.Range("E14").Value = "  +2.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6230"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").Value = "30.354.36"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "243.92"
$ws.Range("E17").Value = "  +8.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  +3.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007295"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.910"
$ws.Range("E22").Value = "  +2.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.893"
$ws.Range("E23").Value = "  +1.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.77"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.073"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.93"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.864"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1029"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.356"
$ws.Range("E29").Value = "  -0.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.030"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.816"
$ws.Range("E31").Value = "  +1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04833"
$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.123"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6953"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.699"
$ws.Range("E35").Value = "  +0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01898"
$ws.Range("E36").Value = "  +4.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  +2.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.995"
$ws.Range("E38").Value = "  +4.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8733"
$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.28"
$ws.Range("E40").Value = "  +3.31%  "

$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.502"
$ws.Range("E42").Value = "  +1.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4039"
$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.143"
$ws.Range("E44").Value = "  +4.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.74"
$ws.Range("E45").Value = "  +6.65%  "

$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.62"
$ws.Range("E47").Value = "  +3.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.515"
$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05516"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.344"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3666"
$ws.Range("E51").Value = "  +1.87%  "
